$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 5 new data columns before column D. This pushes the existing
#    quarterly columns D:H (1400/09 .. 1401/09) to I:M, preserving their
#    values/styles, and leaves fresh (unformatted-width) columns D:H ready
#    for the five new, older quarters (1399/06 .. 1400/06).
# ---------------------------------------------------------------------------
$ws.Range("D1:H1").EntireColumn.Insert()

# Column widths for the newly inserted columns (matches the D:H pattern
# used elsewhere in the sheet: 29,29,31,29,29). Read the widths back from
# the (now shifted) old columns I ("29") and J ("31") so the values
# round-trip through the same char<->pixel quantisation Excel applies,
# instead of hard-coding magic literals that get rounded differently.
$w29 = $ws.Columns.Item(9).ColumnWidth
$w31 = $ws.Columns.Item(10).ColumnWidth
$ws.Columns.Item(4).ColumnWidth = $w29
$ws.Columns.Item(5).ColumnWidth = $w29
$ws.Columns.Item(6).ColumnWidth = $w31
$ws.Columns.Item(7).ColumnWidth = $w29
$ws.Columns.Item(8).ColumnWidth = $w29

# ---------------------------------------------------------------------------
# 2. Header row 8: quarter labels for the new columns D:H.
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("E8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H8").Value = "فصل دوم منتهی به 1400/06"

# ---------------------------------------------------------------------------
# 3. Row 9: publish-date labels for the new columns D:H.
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = "1400-09-14 (5)"
$ws.Range("E9").Value = "1400-10-29 (2)"
$ws.Range("F9").Value = "1401-03-15 (9)"
$ws.Range("G9").Value = "1401-04-29 (2)"
$ws.Range("H9").Value = "1401-09-13 (5)"

# ---------------------------------------------------------------------------
# 4. Financial data for rows 11-27, columns D:H (new quarters).
# ---------------------------------------------------------------------------
$ws.Range("D11").Value = 118524
$ws.Range("E11").Value = 124723
$ws.Range("F11").Value = 236553
$ws.Range("G11").Value = 224878
$ws.Range("H11").Value = 192078

$ws.Range("D12").Value = -62912
$ws.Range("E12").Value = -68470
$ws.Range("F12").Value = -136806
$ws.Range("G12").Value = -137368
$ws.Range("H12").Value = -75902

$ws.Range("D13").Value = 55611
$ws.Range("E13").Value = 56253
$ws.Range("F13").Value = 99748
$ws.Range("G13").Value = 87510
$ws.Range("H13").Value = 116176

$ws.Range("D14").Value = -2201
$ws.Range("E14").Value = -2337
$ws.Range("F14").Value = -3412
$ws.Range("G14").Value = -4091
$ws.Range("H14").Value = -4076

$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"

$ws.Range("D16").Value = 32413
$ws.Range("E16").Value = 417
$ws.Range("F16").Value = -9733
$ws.Range("G16").Value = 582
$ws.Range("H16").Value = -7115

$ws.Range("D17").Value = 85823
$ws.Range("E17").Value = 54333
$ws.Range("F17").Value = 86602
$ws.Range("G17").Value = 84000
$ws.Range("H17").Value = 104985

$ws.Range("D18").Value = -2078
$ws.Range("E18").Value = -2567
$ws.Range("F18").Value = -3550
$ws.Range("G18").Value = -3305
$ws.Range("H18").Value = -519

$ws.Range("D19").Value = 622
$ws.Range("E19").Value = 1682
$ws.Range("F19").Value = -4612
$ws.Range("G19").Value = -233
$ws.Range("H19").Value = 43

$ws.Range("D20").Value = 84367
$ws.Range("E20").Value = 53449
$ws.Range("F20").Value = 78441
$ws.Range("G20").Value = 80462
$ws.Range("H20").Value = 104510

$ws.Range("D21").Value = -7226
$ws.Range("E21").Value = -7215
$ws.Range("F21").Value = -3623
$ws.Range("G21").Value = -9271
$ws.Range("H21").Value = -8712

$ws.Range("D22").Value = 77141
$ws.Range("E22").Value = 46234
$ws.Range("F22").Value = 74818
$ws.Range("G22").Value = 71191
$ws.Range("H22").Value = 95798

$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"

$ws.Range("D24").Value = 77141
$ws.Range("E24").Value = 46234
$ws.Range("F24").Value = 74818
$ws.Range("G24").Value = 71191
$ws.Range("H24").Value = 95798

$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

$ws.Range("D26").Value = 113361
$ws.Range("E26").Value = 91993
$ws.Range("F26").Value = 102357
$ws.Range("G26").Value = 107158
$ws.Range("H26").Value = 95781

$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
